$d = $word.ActiveDocument

$rng = $d.Content.Duplicate
$rng.Find.Execute("Course Feeds.", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)

if ($rng.Find.Found) {
    $start = $rng.Start
    # Replace the whole run's text with "Course " first
    $full = $d.Range($start, $start + 13)
    $full.Text = "Course "

    $afterCourse = $d.Range($start + 7, $start + 7)
    $afterCourse.InsertAfter("Instance ")

    $afterInstance = $d.Range($start + 16, $start + 16)
    $afterInstance.InsertAfter("Feeds.")
}
